$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2: tag existing participant 1 with a DataType ---
$ws.Range("G2").Value = "EDA"

# --- Row 3: new participant 2, same session window as participant 1 ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "12-Jan-2019 11:48:00"
$ws.Range("C3").Value = "12-Jan-2019 12:26:00"
$ws.Range("D3").Formula = "=(HOUR(`$C3-`$B3)*3600)+(MINUTE(`$C3-`$B3)*60)+(SECOND(`$C3-`$B3))"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = "In"
$ws.Range("G3").Value = "Indoor"

# --- Row 4: new participant 3, separate Strava session ---
$ws.Range("A4").Value = 3
$ws.Range("G4").Value = "Strava"
$ws.Range("B4").Value = "16-Aug-2019 10:33:00"
$ws.Range("C4").Value = "16-Aug-2019 10:50:05"
$ws.Range("D4").Formula = "=(HOUR(`$C4-`$B4)*3600)+(MINUTE(`$C4-`$B4)*60)+(SECOND(`$C4-`$B4))"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = "In"

# --- Row 1: new column header (typed last) ---
$ws.Range("G1").Value = "DataType"

# --- Selection left on E7:E8 ---
$ws.Range("E7:E8").Select()
